$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.180.62"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "2.223.01"
$ws.Range("E3").Value = "  +0.72%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "293.88"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.94%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "87.78"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +0.60%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "30.62"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.68%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "50.78"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +6.54%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0781"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.64%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.114"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +3.50%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.42"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "2.567.65"
$ws.Range("E15").Value = "  +0.68%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "13.82"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").Value = "2.215.98"
$ws.Range("E17").Value = "  +0.76%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.737"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("D19").Value = "40.104.22"
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("D20").Value = "0.0₃0891"
$ws.Range("E20").Value = "  +1.19%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "11.24"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.95%  "
$ws.Range("E22").Value = "  -0.10%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "65.63"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.39%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "236.17"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("E26").Value = "  +1.66%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.82"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").Value = "  +3.31%  "
$ws.Range("E29").Value = "  +1.64%  "
$ws.Range("E30").Value = "  -5.83%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "158.87"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +4.20%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "31.85"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("E33").Value = "  -0.03%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.96"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("E35").Value = "  +7.45%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.0714"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.32%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.34"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.26%  "
$ws.Range("E38").Value = "  +1.78%  "
$ws.Range("E39").Value = "  +4.31%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.0995"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.86%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "15.68"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "2.081.63"
$ws.Range("E42").Value = "  -0.55%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "3.76"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.30%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "19.23"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +10.55%  "
$ws.Range("E45").Value = "  +1.36%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "10.05"
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.76"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +4.20%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.92"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -10.55%  "
$ws.Range("D49").Value = "2.440.51"
$ws.Range("E49").Value = "  +0.58%  "
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.12"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +4.36%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.48"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +2.22%  "
